$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the tiny floating point adjustment on row 17, column A
$ws.Range("A17").Value = 45874.66686496528

# Add new row 18 with the latest reading
$ws.Range("A18").Value = 45874.70852650376
$ws.Range("B18").Value = 2025
$ws.Range("C18").Value = 19
$ws.Range("D18").Value = 21.17
$ws.Range("E18").Value = 72.84999999999999
$ws.Range("F18").Value = 156.45
$ws.Range("G18").Value = 9.220000000000001
$ws.Range("H18").Value = "ESE"
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = "17:00:16"

# Match style of A column (date/time format) used by previous rows
$ws.Range("A18").NumberFormat = $ws.Range("A17").NumberFormat
